$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.45
$ws.Range("D5").Value = 1.32
$ws.Range("F5").Value = 1.03
$ws.Range("E6").Value = 1.33
